$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-10-02 Thursday" "2025-10-03 Friday"

Replace-Text "23×46=1058" "27×46=1242"
Replace-Text "29×18=522" "60×82=4920"
Replace-Text "78×71=5538" "76×47=3572"
Replace-Text "17×26=442" "36×94=3384"
Replace-Text "40×52=2080" "58×73=4234"

Replace-Text "69×46=3174" "49×57=2793"
Replace-Text "41×51=2091" "98×82=8036"
Replace-Text "90×65=5850" "67×91=6097"
Replace-Text "79×67=5293" "50×24=1200"
Replace-Text "54×50=2700" "82×14=1148"

Replace-Text "50×18=900" "85×17=1445"
Replace-Text "46×96=4416" "83×63=5229"
Replace-Text "54×24=1296" "78×13=1014"
Replace-Text "98×22=2156" "29×21=609"
Replace-Text "89×16=1424" "96×88=8448"

Replace-Text "13×97=1261" "95×55=5225"
Replace-Text "96×52=4992" "22×39=858"
Replace-Text "70×94=6580" "17×66=1122"
Replace-Text "49×81=3969" "92×79=7268"
Replace-Text "86×28=2408" "44×70=3080"

Replace-Text "72×15=1080" "25×32=800"
Replace-Text "59×35=2065" "94×74=6956"
Replace-Text "82×94=7708" "68×99=6732"
Replace-Text "21×50=1050" "98×85=8330"
Replace-Text "15×12=180" "35×14=490"
